# Applies the "Status code check 2.0" update: appends newly-crawled URL/status rows
# to both the "Canada FR" (sheet1) and "Canada EN" (sheet2) worksheets.

$wb = $excel.ActiveWorkbook

$wsFR = $wb.Worksheets.Item("Canada FR")
$wsEN = $wb.Worksheets.Item("Canada EN")

# Canada FR (sheet1): add rows 37-71
$wsFR.Cells.Item(37, 1).Value = 'https://www.institutabbvie.ca/Bienvenue.html'
$wsFR.Cells.Item(37, 2).Value = 200
$wsFR.Cells.Item(38, 1).Value = 'https://www.institutabbvie.ca/etc.clientlibs/abbvie-pro/components/content/ca-meta-navigation/clientlibs.min.css'
$wsFR.Cells.Item(38, 2).Value = 200
$wsFR.Cells.Item(39, 1).Value = 'https://www.institutabbvie.ca/etc.clientlibs/abbvie-pro/clientlibs/assets/resources/images/icons/Mobile_nav.png'
$wsFR.Cells.Item(39, 2).Value = 200
$wsFR.Cells.Item(40, 1).Value = 'https://www.institutabbvie.ca/cdn-cgi/scripts/5c5dd728/cloudflare-static/email-decode.min.js'
$wsFR.Cells.Item(40, 2).Value = 200
$wsFR.Cells.Item(41, 1).Value = 'https://www.institutabbvie.ca/etc.clientlibs/abbvie-pro/clientlibs/allergan-pro-ami/publish-header.min.css'
$wsFR.Cells.Item(41, 2).Value = 200
$wsFR.Cells.Item(42, 1).Value = 'https://www.institutabbvie.ca/etc.clientlibs/abbvie-pro/clientlibs/allergan-pro-ami/publish-header.min.js'
$wsFR.Cells.Item(42, 2).Value = 200
$wsFR.Cells.Item(43, 1).Value = 'https://www.institutabbvie.ca/etc.clientlibs/abbvie-pro/clientlibs/assets/resources/images/icons/mobile-menu-close.png'
$wsFR.Cells.Item(43, 2).Value = 200
$wsFR.Cells.Item(44, 1).Value = 'https://www.institutabbvie.ca/etc.clientlibs/abbvie-pro/components/content/button-link/clientlibs.min.css'
$wsFR.Cells.Item(44, 2).Value = 200
$wsFR.Cells.Item(45, 1).Value = 'https://www.institutabbvie.ca/etc.clientlibs/abbvie-pro/clientlibs/publish-footer.min.js'
$wsFR.Cells.Item(45, 2).Value = 200
$wsFR.Cells.Item(46, 1).Value = 'https://www.institutabbvie.ca/etc.clientlibs/abbvie-pro/components/content/button-redirect-logout/clientlibs.min.js'
$wsFR.Cells.Item(46, 2).Value = 200
$wsFR.Cells.Item(47, 1).Value = 'https://www.institutabbvie.ca/content/dam/allergan-pro-ami/ca/L_AbbVie%20Institute%20RGB%20F-01-2.png/_jcr_content/renditions/cq5dam.web.1280.1280.png'
$wsFR.Cells.Item(47, 2).Value = 200
$wsFR.Cells.Item(48, 1).Value = 'https://www.institutabbvie.ca/etc.clientlibs/abbvie-pro/components/content/button-link/clientlibs.min.js'
$wsFR.Cells.Item(48, 2).Value = 200
$wsFR.Cells.Item(49, 1).Value = 'https://www.institutabbvie.ca/etc.clientlibs/abbvie-pro/components/content/headline-text/clientlibs.min.js'
$wsFR.Cells.Item(49, 2).Value = 200
$wsFR.Cells.Item(50, 1).Value = 'https://www.institutabbvie.ca/etc.clientlibs/abbvie-pro/components/content/ca-meta-navigation/clientlibs.min.js'
$wsFR.Cells.Item(50, 2).Value = 200
$wsFR.Cells.Item(51, 1).Value = 'https://www.institutabbvie.ca/etc.clientlibs/abbvie-pro/clientlibs/assets/resources/fonts/hcpicon/hcpicon.ttf?q0neb3'
$wsFR.Cells.Item(51, 2).Value = 200
$wsFR.Cells.Item(52, 1).Value = 'https://www.institutabbvie.ca/etc.clientlibs/abbvie-pro/clientlibs/assets/resources/fonts/roboto/roboto_regular.woff2'
$wsFR.Cells.Item(52, 2).Value = 200
$wsFR.Cells.Item(53, 1).Value = 'https://www.institutabbvie.ca/etc.clientlibs/abbvie-pro/clientlibs/assets/resources/fonts/roboto/roboto_bold.woff2'
$wsFR.Cells.Item(53, 2).Value = 200
$wsFR.Cells.Item(54, 1).Value = 'https://consent.trustarc.com/v2/notice/v8idiw'
$wsFR.Cells.Item(54, 2).Value = 200
$wsFR.Cells.Item(55, 1).Value = 'https://www.institutabbvie.ca/etc.clientlibs/abbvie-pro/components/content/image-extension/clientlibs.min.js'
$wsFR.Cells.Item(55, 2).Value = 200
$wsFR.Cells.Item(56, 1).Value = 'https://www.institutabbvie.ca/etc.clientlibs/abbvie-pro/components/content/overlay-component/clientlibs.min.js'
$wsFR.Cells.Item(56, 2).Value = 200
$wsFR.Cells.Item(57, 1).Value = 'https://www.institutabbvie.ca/etc.clientlibs/clientlibs/granite/jquery/granite/csrf.min.js'
$wsFR.Cells.Item(57, 2).Value = 200
$wsFR.Cells.Item(58, 1).Value = 'https://consent.trustarc.com/v2/asset/trustarc-logo-xs.svg'
$wsFR.Cells.Item(58, 2).Value = 200
$wsFR.Cells.Item(59, 1).Value = 'https://consent.trustarc.com/v2/asset/ic-close.svg'
$wsFR.Cells.Item(59, 2).Value = 200
$wsFR.Cells.Item(60, 1).Value = 'https://consent.trustarc.com/v2/asset/latin.woff2'
$wsFR.Cells.Item(60, 2).Value = 200
$wsFR.Cells.Item(61, 1).Value = 'https://consent.trustarc.com/v2/asset/11:37:54.164v8idiw_AbbVieID-logo.png'
$wsFR.Cells.Item(61, 2).Value = 200
$wsFR.Cells.Item(62, 1).Value = 'https://www.institutabbvie.ca/libs/granite/csrf/token.json'
$wsFR.Cells.Item(62, 2).Value = 200
$wsFR.Cells.Item(63, 1).Value = 'https://consent.trustarc.com/v2/consentcategories/getnonemptyindexes?cmId=v8idiw&referer=&fullURL=https%3A%2F%2Fwww.institutabbvie.ca%2FBienvenue.html&category='
$wsFR.Cells.Item(63, 2).Value = 200
$wsFR.Cells.Item(64, 1).Value = 'https://consent-reporting.trustarc.com/api/user-action/bannermsg?action=views&domain=v8idiw&behavior=implied&country=bd&language=en&rand=0.6989222192408637&session=95c8b2ea-72c3-45d6-bb47-9c4a244ef874&userType=NEW'
$wsFR.Cells.Item(64, 2).Value = 202
$wsFR.Cells.Item(65, 1).Value = 'https://consent-reporting.trustarc.com/api/user-action/log?action=impression&domain=v8idiw&behavior=implied&country=bd&language=en&rand=0.21912041421705042&session=95c8b2ea-72c3-45d6-bb47-9c4a244ef874&userType=NEW'
$wsFR.Cells.Item(65, 2).Value = 202
$wsFR.Cells.Item(66, 1).Value = 'https://www.institutabbvie.ca/content/dam/allergan-pro-ami/ca/ami_home3.png/_jcr_content/renditions/cq5dam.web.1280.1280.png'
$wsFR.Cells.Item(66, 2).Value = 200
$wsFR.Cells.Item(67, 1).Value = 'https://consent.trustarc.com/v2/asset/ic-error.svg'
$wsFR.Cells.Item(67, 2).Value = 200
$wsFR.Cells.Item(68, 1).Value = 'https://consent.trustarc.com/v2/asset/ic-close-white.svg'
$wsFR.Cells.Item(68, 2).Value = 200
$wsFR.Cells.Item(69, 1).Value = 'https://www.institutabbvie.ca/bin/public/abbvie-commons/hreflangs?resourcePath=/content/allergan-pro-ami/ca/fr/Bienvenue/jcr:content'
$wsFR.Cells.Item(69, 2).Value = 200
$wsFR.Cells.Item(70, 1).Value = 'https://www.institutabbvie.ca/content/allergan-pro-ami/ca/fr/jcr:content/header/ca_header_area/image-extension/item_1.coreimg.png/1663736091306-L_AbbVie%20Institute%20RGB%20F-01-2.png'
$wsFR.Cells.Item(70, 2).Value = 302
$wsFR.Cells.Item(71, 1).Value = 'https://www.institutabbvie.ca/content/allergan-pro-ami/ca/fr/jcr%3acontent/header/ca_header_area/image-extension/item_1.coreimg.png/1663736100272.png'
$wsFR.Cells.Item(71, 2).Value = 200

# Canada EN (sheet2): add rows 37-71
$wsEN.Cells.Item(37, 1).Value = 'https://www.abbvieinstitute.ca/Welcome.html'
$wsEN.Cells.Item(37, 2).Value = 200
$wsEN.Cells.Item(38, 1).Value = 'https://www.abbvieinstitute.ca/etc.clientlibs/abbvie-pro/clientlibs/allergan-pro-ami/publish-header.min.js'
$wsEN.Cells.Item(38, 2).Value = 200
$wsEN.Cells.Item(39, 1).Value = 'https://www.abbvieinstitute.ca/etc.clientlibs/abbvie-pro/clientlibs/assets/resources/images/icons/Mobile_nav.png'
$wsEN.Cells.Item(39, 2).Value = 200
$wsEN.Cells.Item(40, 1).Value = 'https://www.abbvieinstitute.ca/etc.clientlibs/abbvie-pro/clientlibs/allergan-pro-ami/publish-header.min.css'
$wsEN.Cells.Item(40, 2).Value = 200
$wsEN.Cells.Item(41, 1).Value = 'https://www.abbvieinstitute.ca/etc.clientlibs/abbvie-pro/clientlibs/assets/resources/images/icons/mobile-menu-close.png'
$wsEN.Cells.Item(41, 2).Value = 200
$wsEN.Cells.Item(42, 1).Value = 'https://www.abbvieinstitute.ca/cdn-cgi/scripts/5c5dd728/cloudflare-static/email-decode.min.js'
$wsEN.Cells.Item(42, 2).Value = 200
$wsEN.Cells.Item(43, 1).Value = 'https://www.abbvieinstitute.ca/etc.clientlibs/abbvie-pro/components/content/ca-meta-navigation/clientlibs.min.css'
$wsEN.Cells.Item(43, 2).Value = 200
$wsEN.Cells.Item(44, 1).Value = 'https://www.abbvieinstitute.ca/etc.clientlibs/abbvie-pro/components/content/button-link/clientlibs.min.css'
$wsEN.Cells.Item(44, 2).Value = 200
$wsEN.Cells.Item(45, 1).Value = 'https://www.abbvieinstitute.ca/etc.clientlibs/abbvie-pro/components/content/image-extension/clientlibs.min.js'
$wsEN.Cells.Item(45, 2).Value = 200
$wsEN.Cells.Item(46, 1).Value = 'https://www.abbvieinstitute.ca/etc.clientlibs/clientlibs/granite/jquery/granite/csrf.min.js'
$wsEN.Cells.Item(46, 2).Value = 200
$wsEN.Cells.Item(47, 1).Value = 'https://www.abbvieinstitute.ca/etc.clientlibs/abbvie-pro/clientlibs/publish-footer.min.js'
$wsEN.Cells.Item(47, 2).Value = 200
$wsEN.Cells.Item(48, 1).Value = 'https://www.abbvieinstitute.ca/content/dam/allergan-pro-ami/ca/AMIColorfulLogo.png/_jcr_content/renditions/cq5dam.web.1280.1280.png'
$wsEN.Cells.Item(48, 2).Value = 200
$wsEN.Cells.Item(49, 1).Value = 'https://www.abbvieinstitute.ca/etc.clientlibs/abbvie-pro/components/content/headline-text/clientlibs.min.js'
$wsEN.Cells.Item(49, 2).Value = 200
$wsEN.Cells.Item(50, 1).Value = 'https://www.abbvieinstitute.ca/etc.clientlibs/abbvie-pro/clientlibs/assets/resources/fonts/roboto/roboto_bold.woff2'
$wsEN.Cells.Item(50, 2).Value = 200
$wsEN.Cells.Item(51, 1).Value = 'https://www.abbvieinstitute.ca/etc.clientlibs/abbvie-pro/clientlibs/assets/resources/fonts/roboto/roboto_regular.woff2'
$wsEN.Cells.Item(51, 2).Value = 200
$wsEN.Cells.Item(52, 1).Value = 'https://www.abbvieinstitute.ca/etc.clientlibs/abbvie-pro/components/content/ca-meta-navigation/clientlibs.min.js'
$wsEN.Cells.Item(52, 2).Value = 200
$wsEN.Cells.Item(53, 1).Value = 'https://www.abbvieinstitute.ca/etc.clientlibs/abbvie-pro/components/content/button-link/clientlibs.min.js'
$wsEN.Cells.Item(53, 2).Value = 200
$wsEN.Cells.Item(54, 1).Value = 'https://www.abbvieinstitute.ca/etc.clientlibs/abbvie-pro/clientlibs/assets/resources/fonts/hcpicon/hcpicon.ttf?q0neb3'
$wsEN.Cells.Item(54, 2).Value = 200
$wsEN.Cells.Item(55, 1).Value = 'https://www.abbvieinstitute.ca/content/dam/allergan-pro-ami/ca/ami_home3.png/_jcr_content/renditions/cq5dam.web.1280.1280.png'
$wsEN.Cells.Item(55, 2).Value = 200
$wsEN.Cells.Item(56, 1).Value = 'https://www.abbvieinstitute.ca/etc.clientlibs/abbvie-pro/components/content/button-redirect-logout/clientlibs.min.js'
$wsEN.Cells.Item(56, 2).Value = 200
$wsEN.Cells.Item(57, 1).Value = 'https://www.abbvieinstitute.ca/etc.clientlibs/abbvie-pro/components/content/overlay-component/clientlibs.min.js'
$wsEN.Cells.Item(57, 2).Value = 200
$wsEN.Cells.Item(58, 1).Value = 'https://consent.trustarc.com/v2/notice/rqwfyo'
$wsEN.Cells.Item(58, 2).Value = 200
$wsEN.Cells.Item(59, 1).Value = 'https://consent.trustarc.com/v2/asset/trustarc-logo-xs.svg'
$wsEN.Cells.Item(59, 2).Value = 200
$wsEN.Cells.Item(60, 1).Value = 'https://consent.trustarc.com/v2/asset/ic-close.svg'
$wsEN.Cells.Item(60, 2).Value = 200
$wsEN.Cells.Item(61, 1).Value = 'https://consent.trustarc.com/v2/asset/latin.woff2'
$wsEN.Cells.Item(61, 2).Value = 200
$wsEN.Cells.Item(62, 1).Value = 'https://consent.trustarc.com/v2/asset/11:51:11.658rqwfyo_AbbVieID-logo.png'
$wsEN.Cells.Item(62, 2).Value = 200
$wsEN.Cells.Item(63, 1).Value = 'https://www.abbvieinstitute.ca/content/allergan-pro-ami/ca/en/jcr:content/header/ca_header_area/image-extension/item_1.coreimg.png/1659941845267-AMIColorfulLogo.png'
$wsEN.Cells.Item(63, 2).Value = 302
$wsEN.Cells.Item(64, 1).Value = 'https://www.abbvieinstitute.ca/libs/granite/csrf/token.json'
$wsEN.Cells.Item(64, 2).Value = 200
$wsEN.Cells.Item(65, 1).Value = 'https://consent-reporting.trustarc.com/api/user-action/log?action=impression&domain=rqwfyo&behavior=implied&country=bd&language=en&rand=0.9242139847605746&session=99f6175d-378f-4a1f-b213-228fba045861&userType=NEW'
$wsEN.Cells.Item(65, 2).Value = 202
$wsEN.Cells.Item(66, 1).Value = 'https://consent-reporting.trustarc.com/api/user-action/bannermsg?action=views&domain=rqwfyo&behavior=implied&country=bd&language=en&rand=0.14611231212024833&session=99f6175d-378f-4a1f-b213-228fba045861&userType=NEW'
$wsEN.Cells.Item(66, 2).Value = 202
$wsEN.Cells.Item(67, 1).Value = 'https://consent.trustarc.com/v2/asset/ic-error.svg'
$wsEN.Cells.Item(67, 2).Value = 200
$wsEN.Cells.Item(68, 1).Value = 'https://consent.trustarc.com/v2/asset/ic-close-white.svg'
$wsEN.Cells.Item(68, 2).Value = 200
$wsEN.Cells.Item(69, 1).Value = 'https://www.abbvieinstitute.ca/content/allergan-pro-ami/ca/en/jcr%3acontent/header/ca_header_area/image-extension/item_1.coreimg.png/1659952345719.png'
$wsEN.Cells.Item(69, 2).Value = 200
$wsEN.Cells.Item(70, 1).Value = 'https://www.abbvieinstitute.ca/bin/public/abbvie-commons/hreflangs?resourcePath=/content/allergan-pro-ami/ca/en/Welcome/jcr:content'
$wsEN.Cells.Item(70, 2).Value = 200
$wsEN.Cells.Item(71, 1).Value = 'https://consent.trustarc.com/v2/consentcategories/getnonemptyindexes?cmId=rqwfyo&referer=&fullURL=https%3A%2F%2Fwww.abbvieinstitute.ca%2FWelcome.html&category='
$wsEN.Cells.Item(71, 2).Value = 200

